# Applies the cryptos.xlsx update described in the commit "Updated cryptos
# list on Sat Apr 29 09:46:29 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep being stored as plain text (as in the
# source workbook) rather than being auto-coerced to numbers by Excel's
# input parsing. A leading apostrophe forces text entry while leaving the
# cell format as General, matching the original file.

$ws.Range("D2").Value = "'29.461.90"
$ws.Range("E2").Value = "'  +0.63%  "

$ws.Range("D3").Value = "'1.912.09"
$ws.Range("E3").Value = "'  +0.04%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "'  +0.65%  "

$ws.Range("D5").Value = "'325.41"
$ws.Range("E5").Value = "'  +1.19%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "'  +0.45%  "

$ws.Range("D7").Value = "'0.4823"
$ws.Range("E7").Value = "'  +2.09%  "

$ws.Range("D8").Value = "'0.4067"
$ws.Range("E8").Value = "'  +0.02%  "

$ws.Range("D9").Value = "'0.08150"
$ws.Range("E9").Value = "'  +1.38%  "

$ws.Range("D10").Value = "'1.013"
$ws.Range("E10").Value = "'  +1.06%  "

$ws.Range("D11").Value = "'23.43"

$ws.Range("D12").Value = "'1.914.86"
$ws.Range("E12").Value = "'  +0.30%  "

$ws.Range("D13").Value = "'6.012"
$ws.Range("E13").Value = "'  +2.02%  "

$ws.Range("D14").Value = "'7.151"
$ws.Range("E14").Value = "'  +0.33%  "

$ws.Range("D15").Value = "'90.21"
$ws.Range("E15").Value = "'  +0.64%  "

$ws.Range("D16").Value = "'0.06787"
$ws.Range("E16").Value = "'  +2.30%  "

$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "'  +0.62%  "

$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = "'  +0.87%  "

$ws.Range("D19").Value = "'17.69"
$ws.Range("E19").Value = "'  +0.07%  "

$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "'  +0.50%  "

$ws.Range("D21").Value = "'29.490.44"
$ws.Range("E21").Value = "'  +0.65%  "

$ws.Range("D22").Value = "'5.625"
$ws.Range("E22").Value = "'  +2.03%  "

$ws.Range("D23").Value = "'11.75"
$ws.Range("E23").Value = "'  +2.51%  "

$ws.Range("D24").Value = "'2.182"
$ws.Range("E24").Value = "'  -0.69%  "

$ws.Range("D25").Value = "'2.130.61"
$ws.Range("E25").Value = "'  -0.32%  "

$ws.Range("D26").Value = "'155.82"
$ws.Range("E26").Value = "'  +0.18%  "

$ws.Range("D27").Value = "'6.385"
$ws.Range("E27").Value = "'  +6.11%  "

$ws.Range("D28").Value = "'20.02"
$ws.Range("E28").Value = "'  +1.18%  "

$ws.Range("D29").Value = "'2.106"
$ws.Range("E29").Value = "'  +0.08%  "

$ws.Range("D30").Value = "'120.04"
$ws.Range("E30").Value = "'  +2.31%  "

$ws.Range("D31").Value = "'1.022"
$ws.Range("E31").Value = "'  -4.58%  "

$ws.Range("D32").Value = "'0.09534"
$ws.Range("E32").Value = "'  +0.12%  "

$ws.Range("D33").Value = "'5.522"
$ws.Range("E33").Value = "'  +2.52%  "

$ws.Range("D34").Value = "'3.565"
$ws.Range("E34").Value = "'  +0.66%  "

$ws.Range("D35").Value = "'1.389"
$ws.Range("E35").Value = "'  -2.24%  "

$ws.Range("D36").Value = "'0.02268"
$ws.Range("E36").Value = "'  +1.10%  "

$ws.Range("D37").Value = "'0.06103"
$ws.Range("E37").Value = "'  +0.42%  "

$ws.Range("D38").Value = "'1.175"
$ws.Range("E38").Value = "'  +0.11%  "

$ws.Range("D39").Value = "'0.5967"
$ws.Range("E39").Value = "'  +2.02%  "

$ws.Range("D40").Value = "'10.78"
$ws.Range("E40").Value = "'  +6.77%  "

$ws.Range("D41").Value = "'7.986"
$ws.Range("E41").Value = "'  -3.11%  "

$ws.Range("D42").Value = "'0.1853"
$ws.Range("E42").Value = "'  +1.05%  "

$ws.Range("D43").Value = "'1.278"
$ws.Range("E43").Value = "'  +0.44%  "

$ws.Range("D44").Value = "'2.393"
$ws.Range("E44").Value = "'  -4.65%  "

$ws.Range("D45").Value = "'12.57"
$ws.Range("E45").Value = "'  +3.62%  "

$ws.Range("D46").Value = "'0.07611"
$ws.Range("E46").Value = "'  -3.56%  "

$ws.Range("D47").Value = "'0.5574"
$ws.Range("E47").Value = "'  +0.84%  "

$ws.Range("D48").Value = "'1.939"
$ws.Range("E48").Value = "'  +0.92%  "

$ws.Range("D49").Value = "'115.78"
$ws.Range("E49").Value = "'  +2.40%  "

# Rows 50 and 51: Aave and MXToken swap positions, each with refreshed price/volume figures
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.415"
$ws.Range("E50").Value = "'  +2.86%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'72.50"
$ws.Range("E51").Value = "'  +1.83%  "

Write-Host "Edit applied successfully"
